$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 396, shifting existing rows 396-449 down to 397-450.
$ws.Rows(396).Insert()

# Populate the newly inserted row 396 with the new weekly price entry.
$ws.Range("A396").Value = 11
$ws.Range("B396").Value = "Vega Monumental Concepción"
$ws.Range("C396").Value = "Bíobío"
$ws.Range("D396").Value = 45154
$ws.Range("D396").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E396").Value = 8
$ws.Range("F396").Value = 100114013
$ws.Range("G396").Value = "Zanahoria"
$ws.Range("H396").Value = "Sin especificar"
$ws.Range("I396").Value = "Primera"
$ws.Range("J396").Value = 400
$ws.Range("K396").Value = 5000
$ws.Range("L396").Value = 6000
$ws.Range("M396").Value = 5375
$ws.Range("N396").Value = "$/saco 20 kilos"
$ws.Range("O396").Value = "Región de Ñuble"
$ws.Range("P396").Value = 269
$ws.Range("Q396").Value = 20
$ws.Range("R396").Value = "Hortaliza"
